$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E15 to use a formula instead of a static value
$ws.Range("E15").Formula = "=ROUND(300/0.5/2.44+40,0)"

# Move the active cell selection from B24 to F12
$ws.Range("F12").Select()

$wb.Save()
